# typo fix: "Database 2" -> "Database 3" in slide titles
# (commit message: "I was talking about Database 3 not 2")

$p = $ppt.ActivePresentation

# Slides whose title's very first run begins with "Database 2: " and can be
# fixed with a straight in-place text substitution on that first run only.
$simpleSlides = @{
    2 = "Database 3: ucddb002 128Hz original (+ "
    3 = "Database 3: ucddb002 100Hz (+ "
    4 = "Database 3: ucddb025 128Hz original (+ "
    5 = "Database 3: ucddb025 100Hz (+ "
    8 = "Database 3: 128Hz - all "
}

foreach ($idx in $simpleSlides.Keys) {
    $slide = $p.Slides.Item($idx)
    $title = $slide.Shapes.Item(2)
    $run1 = $title.TextFrame.TextRange.Runs(1)
    $run1.Text = $simpleSlides[$idx]
}

# Slide 9's title run reads "Database 2: 100Hz - all " followed by more runs
# ("data", " ", "sets", " ", "info"). In the target deck only the
# "Database 2: " prefix becomes "Database 3: " and, unlike the other slides,
# that prefix ends up split out into its own separate run (the remainder,
# "100Hz - all ", stays as a second run). Replacing just the first 12
# characters (the "Database 2: " prefix) reproduces that same run split.
$slide9 = $p.Slides.Item(9)
$title9 = $slide9.Shapes.Item(2)
$tr9 = $title9.TextFrame.TextRange
$prefix9 = $tr9.Characters(1, 12)
$prefix9.Text = "Database 3: "
